$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated cell values from the latest crypto data refresh.
# Cells whose new value is a bare number string (e.g. "1.00", "0.168")
# are pre-formatted as Text so Excel keeps the exact text instead of
# silently converting it to a numeric value (which would drop things
# like trailing zeros or switch to scientific notation).
$ws.Range("D2").Value = "61.792.81"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "3.409.22"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "413.05"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.58"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -2.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.725"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("E10").Value = "  -5.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.71"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000216"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.15"
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").Value = "3.941.97"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.41"
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("D17").Value = "3.434.03"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.63"
$ws.Range("E18").Value = "  +3.93%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Value = "61.783.82"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "484.38"
$ws.Range("E21").Value = "  +13.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.65"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.28"
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.07"
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.32"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.77"
$ws.Range("E26").Value = "  +10.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.11"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.76"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.72"
$ws.Range("E29").Value = "  +1.94%  "
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.86"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.112"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.91"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.67"
$ws.Range("E36").Value = "  +8.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0485"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.04"
$ws.Range("E39").Value = "  +4.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.23"
$ws.Range("E40").Value = "  +5.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.325"
$ws.Range("E41").Value = "  +4.17%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.07"
$ws.Range("E44").Value = "  +5.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.57"
$ws.Range("E45").Value = "  +6.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.20"
$ws.Range("E46").Value = "  +2.50%  "
$ws.Range("E47").Value = "  +17.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.40"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.23"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("E50").Value = "  +17.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.07"
$ws.Range("E51").Value = "  +13.23%  "
